$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold numeric-looking text (e.g. "245.10",
# "-0.73%"). Force the cells to Text format first so Excel keeps the literal
# string instead of auto-converting to a number/percentage.
$priceVolCells = @("D2","E2","D3","E3","D4","E4","D5","E5","D6","E6","D7","E7",`
  "D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14",`
  "D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20",`
  "D21","E21","D22","E22","D23","E23","E24","D25","E25","E26","E27","D28","E28",`
  "D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45",`
  "D46","E46","D47","E47","D48","E48","D49","E49","D50","E50")

foreach ($addr in $priceVolCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Rows that only need Price (D) / Volume(1h) (E) value updates
$ws.Range("D2").Value = "245.10"
$ws.Range("E2").Value = "-0.73%"

$ws.Range("D3").Value = "27.22"
$ws.Range("E3").Value = "2.78%"

$ws.Range("D4").Value = "5.106"
$ws.Range("E4").Value = "0.82%"

$ws.Range("D5").Value = "0.05711"
$ws.Range("E5").Value = "2.01%"

$ws.Range("D6").Value = "6.489"
$ws.Range("E6").Value = "0.06%"

$ws.Range("D7").Value = "0.8192"
$ws.Range("E7").Value = "0.74%"

$ws.Range("D8").Value = "0.8557"
$ws.Range("E8").Value = "1.24%"

# Rows 9-20: coin list re-ordered / re-ranked. Update Coin(B), Link(C), Price(D), Volume(E)
$ws.Range("B9").Value = "MandalaExchangeToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D9").Value = "0.06932"
$ws.Range("E9").Value = "-0.96%"

$ws.Range("B10").Value = "BitrueCoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D10").Value = "0.02826"
$ws.Range("E10").Value = "-0.87%"

$ws.Range("B11").Value = "BitMartToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D11").Value = "0.09400"
$ws.Range("E11").Value = "0.08%"

$ws.Range("B12").Value = "BitForexToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D12").Value = "0.001521"
$ws.Range("E12").Value = "0.57%"

$ws.Range("B13").Value = "CoinExToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D13").Value = "0.04028"
$ws.Range("E13").Value = "-13.02%"

$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "0.0005974"
$ws.Range("E14").Value = "0.07%"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.006210"
$ws.Range("E15").Value = "0.83%"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.511"
$ws.Range("E16").Value = "-2.55%"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "3.005"
$ws.Range("E17").Value = "-0.30%"

$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "2.230"
$ws.Range("E18").Value = "8.49%"

$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3165"
$ws.Range("E19").Value = "-1.29%"

$ws.Range("B20").Value = "WazirX"
$ws.Range("C20").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D20").Value = "0.1330"
$ws.Range("E20").Value = "-0.54%"

# Rows 21-28: Price (D) / Volume(1h) (E) updates only
$ws.Range("D21").Value = "0.03204"
$ws.Range("E21").Value = "0.19%"

$ws.Range("D22").Value = "0.1302"
$ws.Range("E22").Value = "0.41%"

$ws.Range("D23").Value = "3.570"
$ws.Range("E23").Value = "-4.56%"

$ws.Range("E24").Value = "1.73%"

$ws.Range("D25").Value = "0.001216"
$ws.Range("E25").Value = "-2.14%"

$ws.Range("E26").Value = "-2.28%"

$ws.Range("E27").Value = "3.11%"

$ws.Range("D28").Value = "0.0001447"
$ws.Range("E28").Value = "3.56%"

# Row 40: Price / Volume update
$ws.Range("D40").Value = "0.03725"
$ws.Range("E40").Value = "1.57%"

# Rows 41-43: coin list re-ordered / re-ranked
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1058"
$ws.Range("E41").Value = "-22.67%"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002459"
$ws.Range("E42").Value = "-7.53%"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.003447"
$ws.Range("E43").Value = "-43.76%"

# Rows 44-50: Price / Volume updates only
$ws.Range("D44").Value = "0.009401"
$ws.Range("E44").Value = "5.06%"

$ws.Range("D45").Value = "0.00005149"
$ws.Range("E45").Value = "-2.50%"

$ws.Range("D46").Value = "0.00000000749"
$ws.Range("E46").Value = "-0.11%"

$ws.Range("D47").Value = "0.1014"
$ws.Range("E47").Value = "-7.82%"

$ws.Range("D48").Value = "0.002509"
$ws.Range("E48").Value = "-4.46%"

$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").Value = "-0.11%"

$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").Value = "-0.11%"
